$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scotland Premiership")

# Match-row data (teams, result, odds) cyclically re-shuffled among grouped rows.
# The 'id' sequence number in column A and Div/Date columns C,D,E stay fixed per row;
# only columns B and F:AC (match id, teams, score, odds) move between rows.

# Row 370
$ws.Range("B370").Value2 = 5169291
$row370 = New-Object 'object[,]' 1,24
$row370[0,0] = 'Livingston'
$row370[0,1] = 'Hearts'
$row370[0,2] = 1
$row370[0,3] = 0
$row370[0,4] = 'H'
$row370[0,5] = 3.5
$row370[0,6] = 3.2
$row370[0,7] = 2.15
$row370[0,8] = 3.4
$row370[0,9] = 3.4
$row370[0,10] = 2.2
$row370[0,11] = 0.25
$row370[0,12] = 1.95
$row370[0,13] = 1.9
$row370[0,14] = 2.5
$row370[0,15] = 1.925
$row370[0,16] = 1.925
$row370[0,17] = 2.4
$row370[0,18] = -1
$row370[0,19] = -1
$row370[0,20] = 0.95
$row370[0,21] = -1
$row370[0,22] = -1
$row370[0,23] = 0.925
$ws.Range("F370:AC370").Value2 = $row370

# Row 371
$ws.Range("B371").Value2 = 5542385
$row371 = New-Object 'object[,]' 1,24
$row371[0,0] = 'Hibernian'
$row371[0,1] = 'Kilmarnock'
$row371[0,2] = 1
$row371[0,3] = 0
$row371[0,4] = 'H'
$row371[0,5] = 1.533
$row371[0,6] = 3.8
$row371[0,7] = 6.5
$row371[0,8] = 1.65
$row371[0,9] = 3.75
$row371[0,10] = 5.75
$row371[0,11] = -0.75
$row371[0,12] = 1.825
$row371[0,13] = 2.025
$row371[0,14] = 2.25
$row371[0,15] = 1.85
$row371[0,16] = 2
$row371[0,17] = 0.6499999999999999
$row371[0,18] = -1
$row371[0,19] = -1
$row371[0,20] = 0.4125
$row371[0,21] = -0.5
$row371[0,22] = -1
$row371[0,23] = 1
$ws.Range("F371:AC371").Value2 = $row371

# Row 372
$ws.Range("B372").Value2 = 5169852
$row372 = New-Object 'object[,]' 1,24
$row372[0,0] = 'Motherwell'
$row372[0,1] = 'Dundee Utd'
$row372[0,2] = 0
$row372[0,3] = 0
$row372[0,4] = 'D'
$row372[0,5] = 2.05
$row372[0,6] = 3.2
$row372[0,7] = 3.75
$row372[0,8] = 1.95
$row372[0,9] = 3.3
$row372[0,10] = 4.2
$row372[0,11] = -0.5
$row372[0,12] = 1.95
$row372[0,13] = 1.9
$row372[0,14] = 2.25
$row372[0,15] = 1.925
$row372[0,16] = 1.925
$row372[0,17] = -1
$row372[0,18] = 2.3
$row372[0,19] = -1
$row372[0,20] = -1
$row372[0,21] = 0.8999999999999999
$row372[0,22] = -1
$row372[0,23] = 0.925
$ws.Range("F372:AC372").Value2 = $row372

# Row 373
$ws.Range("B373").Value2 = 5169604
$row373 = New-Object 'object[,]' 1,24
$row373[0,0] = 'Ross County'
$row373[0,1] = 'Aberdeen'
$row373[0,2] = 1
$row373[0,3] = 1
$row373[0,4] = 'D'
$row373[0,5] = 3.6
$row373[0,6] = 3.2
$row373[0,7] = 2.1
$row373[0,8] = 3.6
$row373[0,9] = 3.2
$row373[0,10] = 2.15
$row373[0,11] = 0.25
$row373[0,12] = 1.95
$row373[0,13] = 1.9
$row373[0,14] = 2.25
$row373[0,15] = 2
$row373[0,16] = 1.85
$row373[0,17] = -1
$row373[0,18] = 2.2
$row373[0,19] = -1
$row373[0,20] = 0.475
$row373[0,21] = -0.5
$row373[0,22] = -0.5
$row373[0,23] = 0.425
$ws.Range("F373:AC373").Value2 = $row373

# Row 374
$ws.Range("B374").Value2 = 5169292
$row374 = New-Object 'object[,]' 1,24
$row374[0,0] = 'St Johnstone'
$row374[0,1] = 'St Mirren'
$row374[0,2] = 3
$row374[0,3] = 0
$row374[0,4] = 'H'
$row374[0,5] = 2.5
$row374[0,6] = 3
$row374[0,7] = 3
$row374[0,8] = 2.9
$row374[0,9] = 3
$row374[0,10] = 2.625
$row374[0,11] = 0
$row374[0,12] = 2.025
$row374[0,13] = 1.825
$row374[0,14] = 2
$row374[0,15] = 1.975
$row374[0,16] = 1.875
$row374[0,17] = 1.9
$row374[0,18] = -1
$row374[0,19] = -1
$row374[0,20] = 1.025
$row374[0,21] = -1
$row374[0,22] = 0.9750000000000001
$row374[0,23] = -1
$ws.Range("F374:AC374").Value2 = $row374

# Row 408
$ws.Range("B408").Value2 = 5169859
$row408 = New-Object 'object[,]' 1,24
$row408[0,0] = 'Dundee Utd'
$row408[0,1] = 'Motherwell'
$row408[0,2] = 0
$row408[0,3] = 1
$row408[0,4] = 'A'
$row408[0,5] = 2.625
$row408[0,6] = 3.2
$row408[0,7] = 2.7
$row408[0,8] = 3
$row408[0,9] = 3.2
$row408[0,10] = 2.5
$row408[0,11] = 0
$row408[0,12] = 2.1
$row408[0,13] = 1.775
$row408[0,14] = 2.25
$row408[0,15] = 1.9
$row408[0,16] = 1.95
$row408[0,17] = -1
$row408[0,18] = -1
$row408[0,19] = 1.5
$row408[0,20] = -1
$row408[0,21] = 0.7749999999999999
$row408[0,22] = -1
$row408[0,23] = 0.95
$ws.Range("F408:AC408").Value2 = $row408

# Row 409
$ws.Range("B409").Value2 = 5169311
$row409 = New-Object 'object[,]' 1,24
$row409[0,0] = 'St Johnstone'
$row409[0,1] = 'Kilmarnock'
$row409[0,2] = 1
$row409[0,3] = 0
$row409[0,4] = 'H'
$row409[0,5] = 2.45
$row409[0,6] = 3
$row409[0,7] = 3.1
$row409[0,8] = 2.4
$row409[0,9] = 2.9
$row409[0,10] = 3.4
$row409[0,11] = -0.25
$row409[0,12] = 2
$row409[0,13] = 1.85
$row409[0,14] = 2
$row409[0,15] = 2
$row409[0,16] = 1.85
$row409[0,17] = 1.4
$row409[0,18] = -1
$row409[0,19] = -1
$row409[0,20] = 1
$row409[0,21] = -1
$row409[0,22] = -1
$row409[0,23] = 0.8500000000000001
$ws.Range("F409:AC409").Value2 = $row409

# Row 410
$ws.Range("B410").Value2 = 5169619
$row410 = New-Object 'object[,]' 1,24
$row410[0,0] = 'Rangers'
$row410[0,1] = 'Aberdeen'
$row410[0,2] = 4
$row410[0,3] = 1
$row410[0,4] = 'H'
$row410[0,5] = 1.333
$row410[0,6] = 5.5
$row410[0,7] = 8
$row410[0,8] = 1.4
$row410[0,9] = 5.25
$row410[0,10] = 6.5
$row410[0,11] = -1.25
$row410[0,12] = 1.825
$row410[0,13] = 2.025
$row410[0,14] = 3
$row410[0,15] = 1.95
$row410[0,16] = 1.9
$row410[0,17] = 0.3999999999999999
$row410[0,18] = -1
$row410[0,19] = -1
$row410[0,20] = 0.825
$row410[0,21] = -1
$row410[0,22] = 0.95
$row410[0,23] = -1
$ws.Range("F410:AC410").Value2 = $row410

# Row 420
$ws.Range("B420").Value2 = 5169623
$row420 = New-Object 'object[,]' 1,24
$row420[0,0] = 'Hibernian'
$row420[0,1] = 'Ross County'
$row420[0,2] = 0
$row420[0,3] = 2
$row420[0,4] = 'A'
$row420[0,5] = 1.6
$row420[0,6] = 3.8
$row420[0,7] = 5.75
$row420[0,8] = 1.4
$row420[0,9] = 4.5
$row420[0,10] = 8
$row420[0,11] = -1.25
$row420[0,12] = 1.95
$row420[0,13] = 1.9
$row420[0,14] = 2.5
$row420[0,15] = 1.95
$row420[0,16] = 1.9
$row420[0,17] = -1
$row420[0,18] = -1
$row420[0,19] = 7
$row420[0,20] = -1
$row420[0,21] = 0.8999999999999999
$row420[0,22] = -1
$row420[0,23] = 0.8999999999999999
$ws.Range("F420:AC420").Value2 = $row420

# Row 421
$ws.Range("B421").Value2 = 5169315
$row421 = New-Object 'object[,]' 1,24
$row421[0,0] = 'Livingston'
$row421[0,1] = 'Aberdeen'
$row421[0,2] = 2
$row421[0,3] = 1
$row421[0,4] = 'H'
$row421[0,5] = 2.625
$row421[0,6] = 3.2
$row421[0,7] = 2.625
$row421[0,8] = 2.45
$row421[0,9] = 3.3
$row421[0,10] = 2.875
$row421[0,11] = 0
$row421[0,12] = 1.775
$row421[0,13] = 2.1
$row421[0,14] = 2.25
$row421[0,15] = 1.9
$row421[0,16] = 1.95
$row421[0,17] = 1.45
$row421[0,18] = -1
$row421[0,19] = -1
$row421[0,20] = 0.7749999999999999
$row421[0,21] = -1
$row421[0,22] = 0.8999999999999999
$row421[0,23] = -1
$ws.Range("F421:AC421").Value2 = $row421

# Row 422
$ws.Range("B422").Value2 = 5169863
$row422 = New-Object 'object[,]' 1,24
$row422[0,0] = 'Rangers'
$row422[0,1] = 'Hearts'
$row422[0,2] = 1
$row422[0,3] = 0
$row422[0,4] = 'H'
$row422[0,5] = 1.285
$row422[0,6] = 5.5
$row422[0,7] = 9.5
$row422[0,8] = 1.2
$row422[0,9] = 7.5
$row422[0,10] = 11
$row422[0,11] = -2
$row422[0,12] = 1.85
$row422[0,13] = 2
$row422[0,14] = 3.5
$row422[0,15] = 1.9
$row422[0,16] = 1.95
$row422[0,17] = 0.2
$row422[0,18] = -1
$row422[0,19] = -1
$row422[0,20] = -1
$row422[0,21] = 1
$row422[0,22] = -1
$row422[0,23] = 0.95
$ws.Range("F422:AC422").Value2 = $row422

# Row 423
$ws.Range("B423").Value2 = 5169862
$row423 = New-Object 'object[,]' 1,24
$row423[0,0] = 'Motherwell'
$row423[0,1] = 'Celtic'
$row423[0,2] = 1
$row423[0,3] = 2
$row423[0,4] = 'A'
$row423[0,5] = 9
$row423[0,6] = 5
$row423[0,7] = 1.333
$row423[0,8] = 9
$row423[0,9] = 7
$row423[0,10] = 1.25
$row423[0,11] = 1.75
$row423[0,12] = 2
$row423[0,13] = 1.85
$row423[0,14] = 3.25
$row423[0,15] = 1.9
$row423[0,16] = 1.95
$row423[0,17] = -1
$row423[0,18] = -1
$row423[0,19] = 0.25
$row423[0,20] = 1
$row423[0,21] = -1
$row423[0,22] = -0.5
$row423[0,23] = 0.475
$ws.Range("F423:AC423").Value2 = $row423

# Row 427
$ws.Range("B427").Value2 = 5169317
$row427 = New-Object 'object[,]' 1,24
$row427[0,0] = 'Hearts'
$row427[0,1] = 'Livingston'
$row427[0,2] = 1
$row427[0,3] = 1
$row427[0,4] = 'D'
$row427[0,5] = 2.05
$row427[0,6] = 3.25
$row427[0,7] = 3.6
$row427[0,8] = 1.8
$row427[0,9] = 3.75
$row427[0,10] = 4.5
$row427[0,11] = -0.75
$row427[0,12] = 2.025
$row427[0,13] = 1.825
$row427[0,14] = 2.5
$row427[0,15] = 2
$row427[0,16] = 1.85
$row427[0,17] = -1
$row427[0,18] = 2.75
$row427[0,19] = -1
$row427[0,20] = -1
$row427[0,21] = 0.825
$row427[0,22] = -1
$row427[0,23] = 0.8500000000000001
$ws.Range("F427:AC427").Value2 = $row427

# Row 429
$ws.Range("B429").Value2 = 5169625
$row429 = New-Object 'object[,]' 1,24
$row429[0,0] = 'Celtic'
$row429[0,1] = 'Ross County'
$row429[0,2] = 2
$row429[0,3] = 1
$row429[0,4] = 'H'
$row429[0,5] = 1.09
$row429[0,6] = 9
$row429[0,7] = 23
$row429[0,8] = 1.071
$row429[0,9] = 15
$row429[0,10] = 21
$row429[0,11] = -3
$row429[0,12] = 1.975
$row429[0,13] = 1.875
$row429[0,14] = 4
$row429[0,15] = 1.9
$row429[0,16] = 1.95
$row429[0,17] = 0.07099999999999995
$row429[0,18] = -1
$row429[0,19] = -1
$row429[0,20] = -1
$row429[0,21] = 0.875
$row429[0,22] = -1
$row429[0,23] = 0.95
$ws.Range("F429:AC429").Value2 = $row429

# Row 447
$ws.Range("B447").Value2 = 5169327
$row447 = New-Object 'object[,]' 1,24
$row447[0,0] = 'St Johnstone'
$row447[0,1] = 'Hearts'
$row447[0,2] = 2
$row447[0,3] = 3
$row447[0,4] = 'A'
$row447[0,5] = 2.7
$row447[0,6] = 3.25
$row447[0,7] = 2.6
$row447[0,8] = 3.1
$row447[0,9] = 3.6
$row447[0,10] = 2.2
$row447[0,11] = 0.25
$row447[0,12] = 1.95
$row447[0,13] = 1.9
$row447[0,14] = 2.5
$row447[0,15] = 1.925
$row447[0,16] = 1.925
$row447[0,17] = -1
$row447[0,18] = -1
$row447[0,19] = 1.2
$row447[0,20] = -1
$row447[0,21] = 0.8999999999999999
$row447[0,22] = 0.925
$row447[0,23] = -1
$ws.Range("F447:AC447").Value2 = $row447

# Row 448
$ws.Range("B448").Value2 = 5169632
$row448 = New-Object 'object[,]' 1,24
$row448[0,0] = 'Dundee Utd'
$row448[0,1] = 'Ross County'
$row448[0,2] = 3
$row448[0,3] = 0
$row448[0,4] = 'H'
$row448[0,5] = 2.6
$row448[0,6] = 3.2
$row448[0,7] = 2.75
$row448[0,8] = 1.909
$row448[0,9] = 3.4
$row448[0,10] = 4.2
$row448[0,11] = -0.5
$row448[0,12] = 1.925
$row448[0,13] = 1.925
$row448[0,14] = 2.25
$row448[0,15] = 1.95
$row448[0,16] = 1.9
$row448[0,17] = 0.909
$row448[0,18] = -1
$row448[0,19] = -1
$row448[0,20] = 0.925
$row448[0,21] = -1
$row448[0,22] = 0.95
$row448[0,23] = -1
$ws.Range("F448:AC448").Value2 = $row448

# Row 466
$ws.Range("B466").Value2 = 5169639
$row466 = New-Object 'object[,]' 1,24
$row466[0,0] = 'Hearts'
$row466[0,1] = 'Aberdeen'
$row466[0,2] = 5
$row466[0,3] = 0
$row466[0,4] = 'H'
$row466[0,5] = 2.375
$row466[0,6] = 3.25
$row466[0,7] = 2.9
$row466[0,8] = 2.1
$row466[0,9] = 3.5
$row466[0,10] = 3.5
$row466[0,11] = -0.25
$row466[0,12] = 1.8
$row466[0,13] = 2.05
$row466[0,14] = 2.75
$row466[0,15] = 1.925
$row466[0,16] = 1.925
$row466[0,17] = 1.1
$row466[0,18] = -1
$row466[0,19] = -1
$row466[0,20] = 0.8
$row466[0,21] = -1
$row466[0,22] = 0.925
$row466[0,23] = -1
$ws.Range("F466:AC466").Value2 = $row466

# Row 467
$ws.Range("B467").Value2 = 5169334
$row467 = New-Object 'object[,]' 1,24
$row467[0,0] = 'Celtic'
$row467[0,1] = 'St Mirren'
$row467[0,2] = 4
$row467[0,3] = 0
$row467[0,4] = 'H'
$row467[0,5] = 1.111
$row467[0,6] = 8
$row467[0,7] = 21
$row467[0,8] = 1.1
$row467[0,9] = 9
$row467[0,10] = 26
$row467[0,11] = -2.25
$row467[0,12] = 1.85
$row467[0,13] = 2
$row467[0,14] = 3.5
$row467[0,15] = 2
$row467[0,16] = 1.85
$row467[0,17] = 0.1000000000000001
$row467[0,18] = -1
$row467[0,19] = -1
$row467[0,20] = 0.8500000000000001
$row467[0,21] = -1
$row467[0,22] = 1
$row467[0,23] = -1
$ws.Range("F467:AC467").Value2 = $row467

# Row 516
$ws.Range("B516").Value2 = 5169887
$row516 = New-Object 'object[,]' 1,24
$row516[0,0] = 'Rangers'
$row516[0,1] = 'Dundee Utd'
$row516[0,2] = 2
$row516[0,3] = 0
$row516[0,4] = 'H'
$row516[0,5] = 1.142
$row516[0,6] = 7
$row516[0,7] = 15
$row516[0,8] = 1.111
$row516[0,9] = 11
$row516[0,10] = 15
$row516[0,11] = -2.5
$row516[0,12] = 1.975
$row516[0,13] = 1.875
$row516[0,14] = 3.5
$row516[0,15] = 1.9
$row516[0,16] = 1.95
$row516[0,17] = 0.111
$row516[0,18] = -1
$row516[0,19] = -1
$row516[0,20] = -1
$row516[0,21] = 0.875
$row516[0,22] = -1
$row516[0,23] = 0.95
$ws.Range("F516:AC516").Value2 = $row516

# Row 517
$ws.Range("B517").Value2 = 5355082
$row517 = New-Object 'object[,]' 1,24
$row517[0,0] = 'St Johnstone'
$row517[0,1] = 'Aberdeen'
$row517[0,2] = 0
$row517[0,3] = 1
$row517[0,4] = 'A'
$row517[0,5] = 2.95
$row517[0,6] = 3.2
$row517[0,7] = 2.4
$row517[0,8] = 3.2
$row517[0,9] = 3.25
$row517[0,10] = 2.25
$row517[0,11] = 0.25
$row517[0,12] = 1.925
$row517[0,13] = 1.925
$row517[0,14] = 2.25
$row517[0,15] = 2.025
$row517[0,16] = 1.825
$row517[0,17] = -1
$row517[0,18] = -1
$row517[0,19] = 1.25
$row517[0,20] = -1
$row517[0,21] = 0.925
$row517[0,22] = -1
$row517[0,23] = 0.825
$ws.Range("F517:AC517").Value2 = $row517

# Row 625
$ws.Range("B625").Value2 = 6844757
$row625 = New-Object 'object[,]' 1,24
$row625[0,0] = 'Motherwell'
$row625[0,1] = 'Aberdeen'
$row625[0,2] = 2
$row625[0,3] = 4
$row625[0,4] = 'A'
$row625[0,5] = 2.625
$row625[0,6] = 3.2
$row625[0,7] = 2.7
$row625[0,8] = 2.45
$row625[0,9] = 3.25
$row625[0,10] = 2.9
$row625[0,11] = 0
$row625[0,12] = 1.8
$row625[0,13] = 2.05
$row625[0,14] = 2.25
$row625[0,15] = 1.9
$row625[0,16] = 1.95
$row625[0,17] = -1
$row625[0,18] = -1
$row625[0,19] = 1.9
$row625[0,20] = -1
$row625[0,21] = 1.05
$row625[0,22] = 0.8999999999999999
$row625[0,23] = -1
$ws.Range("F625:AC625").Value2 = $row625

# Row 626
$ws.Range("B626").Value2 = 6844754
$row626 = New-Object 'object[,]' 1,24
$row626[0,0] = 'Celtic'
$row626[0,1] = 'St Mirren'
$row626[0,2] = 2
$row626[0,3] = 1
$row626[0,4] = 'H'
$row626[0,5] = 1.181
$row626[0,6] = 7
$row626[0,7] = 13
$row626[0,8] = 1.2
$row626[0,9] = 7
$row626[0,10] = 13
$row626[0,11] = -1.75
$row626[0,12] = 1.85
$row626[0,13] = 2
$row626[0,14] = 3
$row626[0,15] = 1.8
$row626[0,16] = 2.05
$row626[0,17] = 0.2
$row626[0,18] = -1
$row626[0,19] = -1
$row626[0,20] = -1
$row626[0,21] = 1
$row626[0,22] = 0
$row626[0,23] = -0
$ws.Range("F626:AC626").Value2 = $row626

# Row 627
$ws.Range("B627").Value2 = 6844758
$row627 = New-Object 'object[,]' 1,24
$row627[0,0] = 'St Johnstone'
$row627[0,1] = 'Kilmarnock'
$row627[0,2] = 2
$row627[0,3] = 1
$row627[0,4] = 'H'
$row627[0,5] = 2.5
$row627[0,6] = 3.1
$row627[0,7] = 2.9
$row627[0,8] = 3.75
$row627[0,9] = 3.1
$row627[0,10] = 2.15
$row627[0,11] = 0.25
$row627[0,12] = 2
$row627[0,13] = 1.85
$row627[0,14] = 2.25
$row627[0,15] = 2
$row627[0,16] = 1.85
$row627[0,17] = 2.75
$row627[0,18] = -1
$row627[0,19] = -1
$row627[0,20] = 1
$row627[0,21] = -1
$row627[0,22] = 1
$row627[0,23] = -1
$ws.Range("F627:AC627").Value2 = $row627

# Row 628
$ws.Range("B628").Value2 = 6844837
$row628 = New-Object 'object[,]' 1,24
$row628[0,0] = 'Hearts'
$row628[0,1] = 'Livingston'
$row628[0,2] = 1
$row628[0,3] = 0
$row628[0,4] = 'H'
$row628[0,5] = 1.727
$row628[0,6] = 3.75
$row628[0,7] = 4.5
$row628[0,8] = 1.45
$row628[0,9] = 4.5
$row628[0,10] = 6.5
$row628[0,11] = -1
$row628[0,12] = 1.825
$row628[0,13] = 2.025
$row628[0,14] = 2.5
$row628[0,15] = 1.9
$row628[0,16] = 1.95
$row628[0,17] = 0.45
$row628[0,18] = -1
$row628[0,19] = -1
$row628[0,20] = 0
$row628[0,21] = -0
$row628[0,22] = -1
$row628[0,23] = 0.95
$ws.Range("F628:AC628").Value2 = $row628

# Row 669
$ws.Range("B669").Value2 = 7451620
$row669 = New-Object 'object[,]' 1,24
$row669[0,0] = 'Aberdeen'
$row669[0,1] = 'Livingston'
$row669[0,2] = 2
$row669[0,3] = 1
$row669[0,4] = 'H'
$row669[0,5] = 1.727
$row669[0,6] = 3.8
$row669[0,7] = 4.5
$row669[0,8] = 1.6
$row669[0,9] = 4
$row669[0,10] = 5.5
$row669[0,11] = -1
$row669[0,12] = 2.05
$row669[0,13] = 1.75
$row669[0,14] = 2.25
$row669[0,15] = 1.8
$row669[0,16] = 2.05
$row669[0,17] = 0.6000000000000001
$row669[0,18] = -1
$row669[0,19] = -1
$row669[0,20] = 0
$row669[0,21] = -0
$row669[0,22] = 0.8
$row669[0,23] = -1
$ws.Range("F669:AC669").Value2 = $row669

# Row 670
$ws.Range("B670").Value2 = 7451733
$row670 = New-Object 'object[,]' 1,24
$row670[0,0] = 'Rangers'
$row670[0,1] = 'St Johnstone'
$row670[0,2] = 2
$row670[0,3] = 0
$row670[0,4] = 'H'
$row670[0,5] = 1.125
$row670[0,6] = 9
$row670[0,7] = 17
$row670[0,8] = 1.125
$row670[0,9] = 8.5
$row670[0,10] = 19
$row670[0,11] = -2.25
$row670[0,12] = 1.925
$row670[0,13] = 1.925
$row670[0,14] = 3.25
$row670[0,15] = 2.025
$row670[0,16] = 1.825
$row670[0,17] = 0.125
$row670[0,18] = -1
$row670[0,19] = -1
$row670[0,20] = -0.5
$row670[0,21] = 0.4625
$row670[0,22] = -1
$row670[0,23] = 0.825
$ws.Range("F670:AC670").Value2 = $row670

# Standalone odds update for upcoming fixture (row 698)
$ws.Range("N698").Value2 = 2.45
$ws.Range("Q698").Value2 = 0
$ws.Range("R698").Value2 = 1.775
$ws.Range("S698").Value2 = 2.1
